$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30: new portfolio data point for 2025-09-14.
# Force the date cell to be stored as literal text (matching the existing
# rows' inline-string dates) rather than letting Excel auto-convert the
# "yyyy-mm-dd" looking text into a date serial number, then strip the
# temporary Text number-format so the cell keeps the workbook's default
# (unstyled) formatting, just like the other data rows.
$ws.Range("A30").NumberFormat = "@"
$ws.Range("A30").Value = "2025-09-14"
$ws.Range("A30").ClearFormats()

$ws.Range("B30").Value = 57.11000061035156
$ws.Range("C30").Value = 715.25
$ws.Range("D30").Value = 321.3999938964844
